$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.254.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.055.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.646"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.80"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.39%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.391"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.66"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +7.92%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.905"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.11"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +18.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.63"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.357.25"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.60"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.055.74"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "37.190.89"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.62%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0894"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.44"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.70"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.52%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.84"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.32"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.38"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.133"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +21.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.123"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.11"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0624"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.58"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.34"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +9.92%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.83"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.96"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +24.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.100"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.90%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.24"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "94.46"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.80"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.378.76"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.93"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.242.50"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.85%  "
